$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "updated on" timestamp in A1
$ws.Range("A1").Value = "更新日期：2025.10.16 13:20:24"

# 2) Insert a new row at 319 (shifts old rows 319-418 down to 320-419)
#    and populate it with the new 圣约送葬人 / GA-3 entry.
$ws.Rows("319:319").Insert()
$ws.Range("A319").Value = "圣约送葬人"
$ws.Range("B319").Value = "GA-3"
$ws.Range("C319").Value = "> 由非助战圣约送葬人累计造成40歼灭数> 3星通关插曲吾导先路GA-3；必须编入非助战圣约送葬人并上场，其他成员仅可编入辅助干员"

# 3) Append a new row after the (now shifted) last row 419, matching the
#    formatting of the preceding row, and populate it with the new
#    真言 / FC-5 entry.
$ws.Rows("420:420").Insert()
$ws.Range("A420").Value = "真言"
$ws.Range("B420").Value = "FC-5"
$ws.Range("C419").Copy()
$ws.Range("C420").PasteSpecial(-4122)
$ws.Range("C420").Value = "> 由非助战真言累计造成60000点元素伤害> 3星通关插曲照我以火FC-5；必须编入非助战真言并上场，且使用真言歼灭至少6名敌人"
